# MTPO - Don1 phase 1 Lua script modification
#
# Duplicate the "Hippo" sheet, place the copy before "Hippo", rename it to
# "Don1-Phase1", and update its stats with the Don1 phase-1 numbers.

$wb = $excel.ActiveWorkbook

$hippo = $wb.Worksheets.Item("Hippo")

# Copy "Hippo" and drop the copy in front of it -> becomes the first sheet.
$hippo.Copy($hippo)
$newSheet = $wb.ActiveSheet
$newSheet.Name = "Don1-Phase1"

# Update the "Perfects / Total / T Error" block for rows 4-6.
$newSheet.Range("A4").Value = 17
$newSheet.Range("B4").Value = 83
$newSheet.Range("C4").Value = 162

$newSheet.Range("A5").ClearContents()
$newSheet.Range("B5").ClearContents()
$newSheet.Range("C5").ClearContents()

$newSheet.Range("A6").ClearContents()
$newSheet.Range("B6").ClearContents()
$newSheet.Range("C6").ClearContents()

# Update the "frames per hit" helper column (J).
$newSheet.Range("J4").Value = 30
$newSheet.Range("J5").Value = 36

# Match the author's last selection on the new sheet.
$newSheet.Range("C5").Select()
